$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 864.6667
$ws.Range("I28").Value = 788
$ws.Range("J28").Value = 903
$ws.Range("K28").Value = 788
$ws.Range("L28").Value = 903
$ws.Range("M28").Value = -303
$ws.Range("N28").Value = -1873

$ws.Range("H32").Value = 9057.727999999999
$ws.Range("I32").Value = 9069.5
$ws.Range("J32").Value = 9055.111000000001
$ws.Range("K32").Value = 9069.5
$ws.Range("L32").Value = 9055.111000000001
$ws.Range("M32").Value = -8743.5
$ws.Range("N32").Value = -9707.111000000001

$ws.Range("H98").Value = 1430
$ws.Range("I98").Value = 986.93335
$ws.Range("J98").Value = 4753
$ws.Range("K98").Value = 986.93335
$ws.Range("L98").Value = 4753
$ws.Range("M98").Value = 511.06665
$ws.Range("N98").Value = -7749

$ws.Range("H113").Value = 2691
$ws.Range("I113").Value = 2587.5
$ws.Range("J113").Value = 3001.5
$ws.Range("K113").Value = 2587.5
$ws.Range("L113").Value = 3001.5
$ws.Range("M113").Value = 666.5
$ws.Range("N113").Value = -9509.5

$ws.Range("H122").Value = 1430
$ws.Range("I122").Value = 986.93335
$ws.Range("J122").Value = 4753
$ws.Range("K122").Value = 2960.80005
$ws.Range("L122").Value = 14259
$ws.Range("M122").Value = -510.8000499999998
$ws.Range("N122").Value = -19159

$ws.Range("H131").Value = 2467.3333
$ws.Range("I131").Value = 1160.8
$ws.Range("J131").Value = 9000
$ws.Range("K131").Value = 3482.4
$ws.Range("L131").Value = 27000
$ws.Range("M131").Value = 1557.6
$ws.Range("N131").Value = -37080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1243.0167
$ws.Range("I32").Value = 1263.2543
$ws.Range("J32").Value = 49
$ws.Range("K32").Value = 1263.2543
$ws.Range("L32").Value = 49
$ws.Range("M32").Value = -976.2543000000001
$ws.Range("N32").Value = -623

$ws.Range("H74").Value = 6615741.5
$ws.Range("I74").Value = 3088735.8
$ws.Range("J74").Value = 27777776
$ws.Range("K74").Value = 3088735.8
$ws.Range("L74").Value = 27777776
$ws.Range("M74").Value = -3087861.8
$ws.Range("N74").Value = -27779524

$ws.Range("H77").Value = 6615741.5
$ws.Range("I77").Value = 3088735.8
$ws.Range("J77").Value = 27777776
$ws.Range("K77").Value = 15443679
$ws.Range("L77").Value = 138888880
$ws.Range("M77").Value = -15439311
$ws.Range("N77").Value = -138897616

$ws.Range("H97").Value = 978.5
$ws.Range("I97").Value = 1021.5238
$ws.Range("J97").Value = 797.8
$ws.Range("K97").Value = 1021.5238
$ws.Range("L97").Value = 797.8
$ws.Range("M97").Value = -525.5238000000001
$ws.Range("N97").Value = -1789.8

$ws.Range("H122").Value = 1086.8889
$ws.Range("I122").Value = 896.3333
$ws.Range("J122").Value = 1468
$ws.Range("K122").Value = 2688.9999
$ws.Range("L122").Value = 4404
$ws.Range("M122").Value = -238.9998999999998
$ws.Range("N122").Value = -9304

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1768.8889
$ws.Range("I105").Value = 1617.1428
$ws.Range("J105").Value = 2300
$ws.Range("K105").Value = 1617.1428
$ws.Range("L105").Value = 2300
$ws.Range("M105").Value = 129.8571999999999
$ws.Range("N105").Value = -5794

$ws.Range("H107").Value = 2036.2
$ws.Range("I107").Value = 656.26666
$ws.Range("J107").Value = 6176
$ws.Range("K107").Value = 656.26666
$ws.Range("L107").Value = 6176
$ws.Range("M107").Value = 1263.73334
$ws.Range("N107").Value = -10016

$ws.Range("H134").Value = 97225144
$ws.Range("I134").Value = 62503750
$ws.Range("J134").Value = 166667920
$ws.Range("K134").Value = 187511250
$ws.Range("L134").Value = 500003760
$ws.Range("M134").Value = -187508715
$ws.Range("N134").Value = -500008830

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3219.6
$ws.Range("I31").Value = 3199
$ws.Range("J31").Value = 3224.75
$ws.Range("K31").Value = 3199
$ws.Range("L31").Value = 3224.75
$ws.Range("M31").Value = -2904
$ws.Range("N31").Value = -3814.75

$ws.Range("H34").Value = 3219.6
$ws.Range("I34").Value = 3199
$ws.Range("J34").Value = 3224.75
$ws.Range("K34").Value = 3199
$ws.Range("L34").Value = 3224.75
$ws.Range("M34").Value = -2997
$ws.Range("N34").Value = -3628.75

$ws.Range("H132").Value = 3871.7632
$ws.Range("I132").Value = 3305.0645
$ws.Range("J132").Value = 6381.4287
$ws.Range("K132").Value = 9915.193499999999
$ws.Range("L132").Value = 19144.2861
$ws.Range("M132").Value = -7385.193499999999
$ws.Range("N132").Value = -24204.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2775.8096
$ws.Range("I122").Value = 3266.2307
$ws.Range("J122").Value = 1978.875
$ws.Range("K122").Value = 9798.6921
$ws.Range("L122").Value = 5936.625
$ws.Range("M122").Value = -7348.6921
$ws.Range("N122").Value = -10836.625

$ws.Range("H126").Value = 9702.846
$ws.Range("I126").Value = 5681.8887
$ws.Range("J126").Value = 18750
$ws.Range("K126").Value = 17045.6661
$ws.Range("L126").Value = 56250
$ws.Range("M126").Value = -14575.6661
$ws.Range("N126").Value = -61190

$ws.Range("H132").Value = 2505.3572
$ws.Range("I132").Value = 2539.1538
$ws.Range("J132").Value = 2066
$ws.Range("K132").Value = 7617.4614
$ws.Range("L132").Value = 6198
$ws.Range("M132").Value = -5087.4614
$ws.Range("N132").Value = -11258

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 49996.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 49996.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 49996.5
$ws.Range("N81").Value = -51992.5

$ws.Range("H84").Value = 49996.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 49996.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 149989.5
$ws.Range("N84").Value = -159973.5

$ws.Range("H122").Value = 3468.0625
$ws.Range("I122").Value = 3210.7144
$ws.Range("J122").Value = 3668.2222
$ws.Range("K122").Value = 9632.143199999999
$ws.Range("L122").Value = 11004.6666
$ws.Range("M122").Value = -7182.143199999999
$ws.Range("N122").Value = -15904.6666

$ws.Range("H132").Value = 2452
$ws.Range("I132").Value = 2176.75
$ws.Range("J132").Value = 3002.5
$ws.Range("K132").Value = 6530.25
$ws.Range("L132").Value = 9007.5
$ws.Range("M132").Value = -4000.25
$ws.Range("N132").Value = -14067.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 907.7
$ws.Range("I126").Value = 888.5
$ws.Range("J126").Value = 984.5
$ws.Range("K126").Value = 2665.5
$ws.Range("L126").Value = 2953.5
$ws.Range("M126").Value = -195.5
$ws.Range("N126").Value = -7893.5

$ws.Range("H132").Value = 2115.818
$ws.Range("I132").Value = 2063.3333
$ws.Range("J132").Value = 2178.8
$ws.Range("K132").Value = 6189.999899999999
$ws.Range("L132").Value = 6536.400000000001
$ws.Range("M132").Value = -3659.999899999999
$ws.Range("N132").Value = -11596.4

$ws.Range("H136").Value = 1342.64
$ws.Range("I136").Value = 1040.6316
$ws.Range("J136").Value = 2299
$ws.Range("K136").Value = 3121.8948
$ws.Range("L136").Value = 6897
$ws.Range("M136").Value = -571.8948
$ws.Range("N136").Value = -11997
